$wb = $excel.ActiveWorkbook
Write-Host ($wb | Get-Member -Name "*Compact*","*Purge*","*Clean*" | Out-String)
$app = $excel
Write-Host ($app | Get-Member -Name "*Compact*","*Purge*","*Clean*" | Out-String)
